# Update 2p4. Added Simscape Multibody Magic Formula Tire for MATLAB R21b and higher.
#
# For every tire sheet, row 6 used to hold an "Inertia" label with three
# plain numbers (10 / 1 / 2) in F6:H6. It is replaced with a new "mjRim"
# row that documents the rim [Mass, Ixx, Iyy] input: a units column (D6),
# a description column (E6), and zeroed F6:H6 values.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Tir_235_50R24", "Tir_213_40R21", "Tir_270_70R22", "Tir_145_70R13", "Tir_430_50R38")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the old trailing blank formatted cells (J6:AA6) that belonged to
    # the "Inertia" row.
    $ws.Range("J6:AA6").ClearContents()

    # New row-6 content.
    $ws.Range("A6").Value = "mjRim"

    $ws.Range("D6").Value = "kg, kg*m^2"
    $ws.Range("D6").Style = "Normal"

    $ws.Range("E6").Value = "Rim [Mass, Ixx, Iyy]"
    $ws.Range("E6").Style = "Normal"

    $ws.Range("F6").Value = 0
    $ws.Range("G6").Value = 0
    $ws.Range("H6").Value = 0

    # Rebuild the "class" highlight conditional formatting so the new D6:H6
    # area isn't covered by the old big rule, and row 6 (A6 / C6:D6) gets its
    # own rule instead (mirrors the authored workbook, which split the single
    # C4:D8/A4:A13-style sqref around the new row).
    $bigRange = $ws.Range("C21:C24,A21:A24,C4:D8,C9:C13,A4:A13")
    $bigCond = $bigRange.FormatConditions
    if ($bigCond.Count -ge 1) {
        $bigCond.Item(1).Delete()
    }

    $splitRanges = @("C21:C24", "A21:A24", "C4:D5", "C9:C13", "A4:A5", "A7:A13", "C7:D8")
    foreach ($r in $splitRanges) {
        $rule = $ws.Range($r).FormatConditions.Add(1, 3, '"class"')
        $rule.Interior.Color = 12937777
    }

    $rowRule = $ws.Range("C6:D6,A6").FormatConditions.Add(1, 3, '"class"')
    $rowRule.Interior.Color = 12937777
}

# View/selection bookkeeping: Excel re-saved with the first sheet active
# (instead of the last), new per-sheet selections, and the last sheet's
# frozen bottom-right pane now highlights the whole new row 6.

$ws2 = $wb.Worksheets.Item("Tir_213_40R21")
$ws2.Activate()
$ws2.Range("K22").Select()

$ws3 = $wb.Worksheets.Item("Tir_270_70R22")
$ws3.Activate()
$ws3.Range("K22").Select()

$ws4 = $wb.Worksheets.Item("Tir_145_70R13")
$ws4.Activate()
$ws4.Range("K22").Select()

$ws5 = $wb.Worksheets.Item("Tir_430_50R38")
$ws5.Activate()
$ws5.Range("A6:XFD6").Select()

$ws1 = $wb.Worksheets.Item("Tir_235_50R24")
$ws1.Activate()
$ws1.Range("E23").Select()
